# mbedXdotShield v1.1.0 BOM update
# - R3 moves from the "1.5K" group to the "220R" group (R3: 1.5K -> 220R)
# - R12, R13 move from the "470R" group to the "1.5K" group (R12,R13: 470R -> 1.5K)
# - R11 stays alone in the "470R" group
# - Qty column bulk-refreshed (matches the reworked schematic / placement pass)
# - Selection moved to F39

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "Parts" (component reference) lists that actually changed content.
# Order matters here only for the underlying shared-string table layout, not
# for correctness, but we set them in the same order the source workbook did.
$ws.Range("E11").Value = "R3, R14, R15, R16, R17, R18, R19, R20, R21, R22, R23, R24, R25, R26"
$ws.Range("E4").Value  = "R2, R5, R8, R12, R13"
$ws.Range("E13").Value = "R11"

# Qty column refresh for every BOM line.
$ws.Range("A4").Value  = 5
$ws.Range("A5").Value  = 3
$ws.Range("A6").Value  = 3
$ws.Range("A7").Value  = 3
$ws.Range("A8").Value  = 3
$ws.Range("A9").Value  = 3
$ws.Range("A10").Value = 3
$ws.Range("A11").Value = 14
$ws.Range("A12").Value = 3
$ws.Range("A13").Value = 1
$ws.Range("A14").Value = 3
$ws.Range("A15").Value = 3
$ws.Range("A16").Value = 3
$ws.Range("A17").Value = 3
$ws.Range("A18").Value = 3
$ws.Range("A19").Value = 3
$ws.Range("A20").Value = 3
$ws.Range("A21").Value = 3
$ws.Range("A22").Value = 3
$ws.Range("A23").Value = 3
$ws.Range("A24").Value = 3
$ws.Range("A25").Value = 3
$ws.Range("A26").Value = 3
$ws.Range("A27").Value = 3
$ws.Range("A28").Value = 3
$ws.Range("A29").Value = 3
$ws.Range("A30").Value = 3
$ws.Range("A31").Value = 3
$ws.Range("A32").Value = 3
$ws.Range("A33").Value = 3

# Match the author's final selection in the saved file.
$ws.Range("F39").Select()
